$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 should share the same (bold/bordered) style as the
# other header cells, e.g. H1. Copy its formatting over before setting values.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2/J2 (plain, unstyled numeric cells, like the rest of row 2)
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
